# Updates the wrss (weighted reel strip symbols) data table: rows are
# reordered (re-sorted) after a "height" value was incorporated into the
# chgSymbols calculation. No columns were added; values per symbol id are
# unchanged, only their row order changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2..23 (columns A..F): symbol, reel1..reel5
$data = @(
    @(301,6,45,30,60,45),
    @(701,3,90,45,97,15),
    @(1203,3,15,15,15,15),
    @(101,9,30,15,60,15),
    @(901,16,15,45,60,60),
    @(201,9,30,15,45,30),
    @(1001,18,30,75,60,72),
    @(501,9,52,30,75,45),
    @(1202,2,10,10,10,10),
    @(601,9,60,67,60,42),
    @(801,3,67,65,52,45),
    @(1201,2,10,10,10,10),
    @(902,1,0,0,0,0),
    @(401,9,48,67,75,45),
    @(802,0,4,5,4,0),
    @(1,0,2,2,2,2),
    @(3,0,3,3,3,3),
    @(502,0,4,0,0,0),
    @(1101,0,15,30,30,0),
    @(2,0,2,2,2,2),
    @(602,0,0,4,0,9),
    @(402,0,0,4,0,0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
